$d = $word.ActiveDocument

# --- Insert "GalloSalvato, " after "Doriana, " (before "IPini") ---
$rng = $d.Content
$found = $rng.Find.Execute("Doriana, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)  # wdCollapseEnd -> collapse to right before "IPini"
$insPos = $rng.Start
$rng.InsertAfter("GalloSalvato, ")
$gs = $d.Range($insPos, $insPos + 12)   # "GalloSalvato" = 12 chars
$gs.Bold = 1
$gs.Italic = 1

# --- Insert ", LupiLupi" after "Linneo" (before " e MakeNao.") ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Linneo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Collapse(0)  # collapse to right after "Linneo"
$insPos2 = $rng2.Start
$rng2.InsertAfter(", LupiLupi")
$comma = $d.Range($insPos2, $insPos2 + 1)       # ","
$comma.Bold = 0
$comma.Italic = 0
$lupi = $d.Range($insPos2 + 1, $insPos2 + 10)   # " LupiLupi"
$lupi.Bold = 1
$lupi.Italic = 1

# --- Move the _GoBack bookmark to sit between the space and "e" of " e MakeNao." ---
$bmPos = $insPos2 + 10 + 1   # skip the inserted text + the following space
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
